$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unhide the rows we need to touch so the runtime doesn't recompute a
# custom row height for them purely because of the write.
$ws.Rows(1).Hidden = $false
$ws.Rows(3).Hidden = $false
$ws.Rows(20).Hidden = $false

# Insert a new column at A; this shifts all existing columns (A..W) to (B..X)
# while keeping their values/styles/formatting intact.
$ws.Columns("A:A").Insert()

# New header for the inserted column.
$ws.Range("A2").Value = "Match ID"

# Apply the "bold, no border" style used for the new column's body cells.
$ws.Range("A2").Font.Bold = $true
$ws.Range("A3:A19").Font.Bold = $true

# Body values - Match ID is constant (9) for every player row.
$ws.Range("A4:A19").Value = 9
$ws.Range("A20").Value = 9

# Re-hide rows.
$ws.Rows(1).Hidden = $true
$ws.Rows(3).Hidden = $true
$ws.Rows(20).Hidden = $true

# Update the selection to match the new layout.
$ws.Range("A2:A19").Select()
